$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet (tab name) from "alpha4F-HW25.xpc" to "alpha4F"
$ws.Name = "alpha4F"

# 2. Append a new row of data (row 16), mirroring the formatting of row 15
#    Column A uses the same style as the rest of column A (apply via CopyPasteSpecial of format)
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.066167906916869
$ws.Range("D16").Value = 0.8099588618061759
$ws.Range("E16").Value = 1.027829301863287
$ws.Range("F16").Value = 1.066167906916869
$ws.Range("G16").Value = 0.8952340248181317
$ws.Range("H16").Value = 1.085230290652888
$ws.Range("I16").Value = 1.040276857054268
$ws.Range("J16").Value = 0.8099588618061759
$ws.Range("K16").Value = 0.9188940818347316
$ws.Range("L16").Value = 0.9925309943758004
$ws.Range("M16").Value = 0.9874495405186033
